$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.142.66"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.540.69"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.18"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.60"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.543.81"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("E10").Value = "  -4.02%  "
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.146.05"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000208"
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.48"
$ws.Range("E15").Value = "  -3.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.542.58"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.322.22"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.03"
$ws.Range("E19").Value = "  +3.43%  "
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("E21").Value = "  -1.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.20"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.602"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.91"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.686.58"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.33"
$ws.Range("E28").Value = "  -4.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.50"
$ws.Range("E32").Value = "  -5.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.161"
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.38"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.534.48"
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.89"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("E39").Value = "  -4.82%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.83"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0863"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("E43").Value = "  -4.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.894"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.88"
$ws.Range("E45").Value = "  -9.55%  "
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.16"
$ws.Range("E47").Value = "  -6.95%  "
$ws.Range("E48").Value = "  -8.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.43"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.20"
$ws.Range("E51").Value = "  -3.83%  "
